$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.482.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = "'2.380.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'503.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = "'130.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E8").Value = '  -2.10%  '
$ws.Range("D9").Value = "'2.389.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.23%  '
$ws.Range("D10").Value = "'0.0988"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = "'4.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.40%  '
$ws.Range("D14").Value = "'2.804.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = "'56.469.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = "'21.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = "'2.339.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = "'10.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.59%  '
$ws.Range("E20").Value = '  -1.65%  '
$ws.Range("D21").Value = "'307.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = "'6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("E23").Value = '  +0.30%  '
$ws.Range("D24").Value = "'65.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = "'0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = "'0.369"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.21%  '
$ws.Range("E27").Value = '  -3.31%  '
$ws.Range("D28").Value = "'7.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.06%  '
$ws.Range("D29").Value = "'172.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.11%  '
$ws.Range("E30").Value = '  -1.62%  '
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  -6.07%  '
$ws.Range("E34").Value = '  -3.42%  '
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("D36").Value = "'17.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("E37").Value = '  -4.97%  '
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").Value = "'0.796"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("E41").Value = '  -4.19%  '
$ws.Range("D42").Value = "'131.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.43%  '
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").Value = "'4.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = "'0.565"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.68%  '
$ws.Range("E46").Value = '  -0.70%  '
$ws.Range("D47").Value = "'242.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.56%  '
$ws.Range("E48").Value = '  -1.79%  '
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Value = "'17.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.27%  '
$ws.Range("E51").Value = '  -1.94%  '
